$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.669.96"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.136.18"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.10"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.33"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.134.76"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.94"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "3.656.82"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "63.535.65"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "3.134.83"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.45"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.91"
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.11"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("E28").Value = "  +5.45%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.90"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.0₃0843"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -6.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.98"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "438.90"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.81"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "2.903.60"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.276"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.89"
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.66"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.21"
$ws.Range("E51").Value = "  -2.17%  "
